$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet by duplicating "2022-Q2" (so that it
#    inherits the exact same header row / column-A styling), then place it
#    immediately before "2022-Q2".
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Row 2: fund 010690
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "010690"
$q3Sheet.Range("C2").Value = "万家互联互通核心资产量化策略混合A"
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "0.55"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "92.41"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "6.86"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0377"
$q3Sheet.Range("H2").Value = 9

# Row 3: fund 010691 (copy row 2's formatting first, then overwrite values)
$q3Sheet.Range("A2:H2").Copy()
$q3Sheet.Range("A3:H3").PasteSpecial(-4122)

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "010691"
$q3Sheet.Range("C3").Value = "万家互联互通核心资产量化策略混合C"
$q3Sheet.Range("D3").Value = "0.15"
$q3Sheet.Range("E3").Value = "92.41"
$q3Sheet.Range("F3").Value = "6.86"
$q3Sheet.Range("G3").Value = "0.0103"
$q3Sheet.Range("H3").Value = 9

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) worksheet: push the three existing data rows
#    down by one and insert the new 2022-Q3 totals at the top.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# extend formatting from row 4 down into the new row 5
$totalSheet.Range("A4:D4").Copy()
$totalSheet.Range("A5:D5").PasteSpecial(-4122)

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q1"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.59

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.62

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.05

# ---------------------------------------------------------------------------
# 3. Restore the originally-selected tab (the last sheet, "2021-Q1").
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
